# Fruta / hortaliza, semanal
# Insert a new weekly record at row 99 (shifting existing rows 99-112 down to 100-113)
# for the "Terminal Hortofrutícola Agro Chillán - Alcachofa" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 99, pushing current rows 99-112 down to 100-113.
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with the new week's data.
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value = "Ñuble"
$ws.Cells.Item(99, 4).Value = 45142
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = 100112013
$ws.Cells.Item(99, 7).Value = "Alcachofa"
$ws.Cells.Item(99, 8).Value = "Argentina(o)"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 50
$ws.Cells.Item(99, 11).Value = 15000
$ws.Cells.Item(99, 12).Value = 15000
$ws.Cells.Item(99, 13).Value = 15000
$ws.Cells.Item(99, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 300
$ws.Cells.Item(99, 17).Value = 50
$ws.Cells.Item(99, 18).Value = "Hortaliza"
